# chore: apply phase2.1.1 release hygiene hardening and refresh README
#
# Applies the OOXML diff across three worksheets:
#   - "00_읽는법" (sheet index 1): refresh the S1/S5 reference-doc wording
#   - "09_출처추적_매트릭스" (sheet index 10): refresh verification timestamps,
#     Notion sync status row, UTF-8 check row, and evidence artifact paths
#   - "10_문서변경이력" (sheet index 11): append the two phase2.1 release rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 00_읽는법 — B3 reference-document summary
# ---------------------------------------------------------------------
$wsReadMe = $wb.Worksheets.Item(1)
$wsReadMe.Range("B3").Value = "S1(20260221 hardening), S2(202603XX full report), S3(Go-Live Gap Closure), S4(spec_sync_report), S5(phase2_1 evidence)"

# ---------------------------------------------------------------------
# 09_출처추적_매트릭스 — verification timestamps + two evidence rows
# ---------------------------------------------------------------------
$wsTrace = $wb.Worksheets.Item(10)

# M-001..M-004 re-verified at 14:17:08Z
$wsTrace.Range("E2").Value = "2026-02-21T14:17:08Z"
$wsTrace.Range("E3").Value = "2026-02-21T14:17:08Z"
$wsTrace.Range("E4").Value = "2026-02-21T14:17:08Z"
$wsTrace.Range("E5").Value = "2026-02-21T14:17:08Z"

# M-005 spec_consistency artifact renamed + re-verified at 14:28:40Z
$wsTrace.Range("C6").Value = "docs/review/mvp_verification_pack/artifacts/phase2_1_pr2_spec_consistency_202603XX.txt"
$wsTrace.Range("E6").Value = "2026-02-21T14:28:40Z"

# M-006, M-007 re-verified at 14:28:40Z
$wsTrace.Range("E7").Value = "2026-02-21T14:28:40Z"
$wsTrace.Range("E8").Value = "2026-02-21T14:28:40Z"

# M-008 Notion sync: now fail-closed/blocked instead of auto-DONE
$wsTrace.Range("B9").Value = "Notion auto-sync blocked state (manual patch operation)"
$wsTrace.Range("C9").Value = "docs/review/mvp_verification_pack/artifacts/phase2_1_pr2_notion_sync_status_202603XX.txt"
$wsTrace.Range("D9").Value = "status=BLOCKED_AUTOMATION"
$wsTrace.Range("E9").Value = "2026-02-21T14:28:40Z"
$wsTrace.Range("F9").Value = "WARN"

# M-009 UTF-8 check wording + artifact rename
$wsTrace.Range("B10").Value = "Phase2.1 UTF-8 validation passed"
$wsTrace.Range("C10").Value = "docs/review/mvp_verification_pack/artifacts/phase2_1_utf8_check_202603XX.txt"
$wsTrace.Range("E10").Value = "2026-02-21T14:28:40Z"

# M-010 re-verified at 14:28:40Z
$wsTrace.Range("E11").Value = "2026-02-21T14:28:40Z"

# ---------------------------------------------------------------------
# 10_문서변경이력 — append phase2.1 and phase2.1-r2 release rows
# ---------------------------------------------------------------------
$wsHistory = $wb.Worksheets.Item(11)

$wsHistory.Range("A5").Value = "v2026.02.21-phase2.1"
$wsHistory.Range("B5").Value = "2026-02-21 23:17:08 +0900"
$wsHistory.Range("C5").Value = "Phase2.1 PR1~PR3 ??(Notion fail-closed, async export-jobs, scheduler self-healing, ??/?? ???)"
$wsHistory.Range("D5").Value = "S1,S2,S3,S4,S5"

$wsHistory.Range("A6").Value = "v2026.02.21-phase2.1-r2"
$wsHistory.Range("B6").Value = "2026-02-21 23:28:40 +09:00"
$wsHistory.Range("C6").Value = "menual consistency refresh (Notion status wording, phase2_1 evidence paths, risk wording normalization)"
$wsHistory.Range("D6").Value = "S1,S2,S3,S4,S5,S6,S7"

Write-Host "phase2.1.1 release hygiene edits applied"
